# Apply crypto price/volume updates per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.427.71"
$ws.Range("E2").Value = "  -1.05%  "
$ws.Range("D3").Value = "2.053.43"
$ws.Range("E3").Value = "  -1.32%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "230.85"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.19%  "
$ws.Range("E6").Value = "  -1.60%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "57.23"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.13%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.388"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.64%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0811"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.30%  "
$ws.Range("E11").Value = "  -2.07%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.76"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.88%  "
$ws.Range("D13").Value = "2.357.32"
$ws.Range("E13").Value = "  -1.26%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.86"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.00%  "
$ws.Range("E15").Value = "  -2.36%  "
$ws.Range("E16").Value = "  -0.87%  "
$ws.Range("D17").Value = "2.044.70"
$ws.Range("E17").Value = "  -1.38%  "
$ws.Range("D18").Value = "37.312.37"
$ws.Range("E18").Value = "  -1.27%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.10"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.56%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "70.02"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.79%  "
$ws.Range("D21").Value = "0.0₃0843"
$ws.Range("E21").Value = "  +0.56%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "227.29"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.29%  "
$ws.Range("E23").Value = "  +0.09%  "
$ws.Range("E24").Value = "  +0.13%  "
$ws.Range("E25").Value = "  -4.39%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.59"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.42%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "167.87"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.38%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.42"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.80%  "
$ws.Range("E29").Value = "  -4.83%  "
$ws.Range("E30").Value = "  -2.45%  "
$ws.Range("E31").Value = "  -2.43%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.56"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.44%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.62"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.47%  "
$ws.Range("E34").Value = "  -2.66%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.43"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.96%  "
$ws.Range("E36").Value = "  +0.15%  "
$ws.Range("E37").Value = "  +0.18%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.26"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.23%  "
$ws.Range("E39").Value = "  -1.56%  "
$ws.Range("E40").Value = "  -5.59%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "17.21"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.44%  "
$ws.Range("D42").Value = "1.493.00"
$ws.Range("E42").Value = "  +2.85%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.89"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.25%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "96.93"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.30%  "
$ws.Range("E46").Value = "  +1.52%  "
$ws.Range("E47").Value = "  -3.77%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.18"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.25%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.92"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.38%  "
$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").Value = "2.244.09"
$ws.Range("E50").Value = "  -1.22%  "
$ws.Range("B51").Value = "FTXToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.70"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -9.58%  "
